# Increase font sizes across the resume document per the mapping:
#   16pt -> 18pt   (name header)
#    9pt -> 10pt   (contact info / dates / bullets / body text)
#   12pt -> 13pt   (section headers)
#   10pt -> 11pt   (overview paragraph)
#   11pt -> 12pt   (job titles)
#
# Every run in this document is sized uniformly per-paragraph, so walk
# every paragraph's Range and bump Font.Size for each distinct size seen.

$d = $word.ActiveDocument

$sizeMap = @{
    16 = 18
    9  = 10
    12 = 13
    10 = 11
    11 = 12
}

foreach ($para in $d.Paragraphs) {
    $pRange = $para.Range
    # Exclude the trailing paragraph-mark character so we don't stamp
    # sz onto the pPr/rPr (paragraph mark run properties) -- only the
    # actual text runs should change size, matching the diff.
    $r = $d.Range($pRange.Start, $pRange.End - 1)
    if ($r.Start -lt $r.End) {
        $cur = $r.Font.Size
        if ($sizeMap.ContainsKey($cur)) {
            $r.Font.Size = $sizeMap[$cur]
        }
    }
}
